# B6-PowerPoint.pptx edit:
#  1. Re-theme the deck from the "Integral" (Red Violet) design to the
#     built-in "Office Theme" colour palette (Design tab -> Office Theme).
#  2. Re-apply the (now complementary) built-in table style to the three
#     tables in the deck (previously a custom "no style, no grid" style).

$p = $ppt.ActivePresentation

# --- 1. Theme colours: Red Violet -> Office -----------------------------
# ThemeColorScheme index order is dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink (MsoThemeColorSchemeIndex 1..12). Values are packed BGR longs,
# i.e. the same encoding VBA's RGB(r,g,b) produces.
$officeColors = @(
    0,          # dk1     000000
    16777215,   # lt1     FFFFFF
    6968388,    # dk2     44546A
    15132391,   # lt2     E7E6E6
    13998939,   # accent1 5B9BD5
    3243501,    # accent2 ED7D31
    10855845,   # accent3 A5A5A5
    49407,      # accent4 FFC000
    12874308,   # accent5 4472C4
    4697456,    # accent6 70AD47
    12673797,   # hlink   0563C1
    7491477     # folHlink 954F72
)

$tcs = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Item($i).RGB = $officeColors[$i - 1]
}

# --- 2. Table styles: custom "Table_0" -> built-in style ----------------
$newStyleId = "{1E6919C4-3617-4327-870D-DE99CFCD9AE6}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
